$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46081 -> 46082, i.e. 2026-02-28 -> 2026-03-01) for every data row (rows 2-299).
for ($r = 2; $r -le 299; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
